# Femacal de La Calera - Papaya: weekly fruit/hortaliza price update.
# Insert two new rows (90:91) for the week of 2023-10-05, shifting the
# existing rows 90-116 down to 92-118.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A90:A91").EntireRow.Insert()

# New row 90: Primera
$ws.Cells.Item(90, 1).Value = 3
$ws.Cells.Item(90, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(90, 3).Value = "Coquimbo"
$ws.Cells.Item(90, 4).Value = 45204
$ws.Cells.Item(90, 5).Value = 5
$ws.Cells.Item(90, 6).Value = "Fruta"
$ws.Cells.Item(90, 7).Value = 100108
$ws.Cells.Item(90, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(90, 9).Value = 100108004
$ws.Cells.Item(90, 10).Value = "Papaya"
$ws.Cells.Item(90, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(90, 12).Value = "Primera"
$ws.Cells.Item(90, 13).Value = 54
$ws.Cells.Item(90, 14).Value = 17000
$ws.Cells.Item(90, 15).Value = 17000
$ws.Cells.Item(90, 16).Value = 17000
$ws.Cells.Item(90, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(90, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(90, 19).Value = 1700
$ws.Cells.Item(90, 20).Value = 10

# New row 91: Segunda
$ws.Cells.Item(91, 1).Value = 3
$ws.Cells.Item(91, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = 45204
$ws.Cells.Item(91, 5).Value = 5
$ws.Cells.Item(91, 6).Value = "Fruta"
$ws.Cells.Item(91, 7).Value = 100108
$ws.Cells.Item(91, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(91, 9).Value = 100108004
$ws.Cells.Item(91, 10).Value = "Papaya"
$ws.Cells.Item(91, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(91, 12).Value = "Segunda"
$ws.Cells.Item(91, 13).Value = 50
$ws.Cells.Item(91, 14).Value = 15000
$ws.Cells.Item(91, 15).Value = 15000
$ws.Cells.Item(91, 16).Value = 15000
$ws.Cells.Item(91, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(91, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(91, 19).Value = 1500
$ws.Cells.Item(91, 20).Value = 10
